$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose updated price text looks like a plain number
# (e.g. "411.94") need to be forced to Text format first, otherwise Excel
# auto-converts the assigned string into a real number (as it does for any
# normal numeric entry) and we would lose the original text-cell semantics.
$numericLookingPriceCells = @(
    "D4", "D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D31", "D33", "D34", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50"
)
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price (column D) updates ---
$ws.Range("D2").Value = '66.378.85'
$ws.Range("D3").Value = '3.545.95'
$ws.Range("D4").Value = '1.01'
$ws.Range("D5").Value = '411.94'
$ws.Range("D6").Value = '127.99'
$ws.Range("D7").Value = '0.638'
$ws.Range("D8").Value = '3.539.12'
$ws.Range("D10").Value = '0.762'
$ws.Range("D11").Value = '0.169'
$ws.Range("D12").Value = '0.0000303'
$ws.Range("D13").Value = '41.63'
$ws.Range("D14").Value = '9.72'
$ws.Range("D15").Value = '4.119.52'
$ws.Range("D17").Value = '19.96'
$ws.Range("D18").Value = '3.563.97'
$ws.Range("D19").Value = '1.10'
$ws.Range("D20").Value = '66.400.75'
$ws.Range("D21").Value = '12.13'
$ws.Range("D22").Value = '440.91'
$ws.Range("D23").Value = '88.24'
$ws.Range("D24").Value = '3.10'
$ws.Range("D25").Value = '12.83'
$ws.Range("D26").Value = '3.37'
$ws.Range("D28").Value = '33.93'
$ws.Range("D31").Value = '12.17'
$ws.Range("D33").Value = '7.19'
$ws.Range("D34").Value = '0.159'
$ws.Range("D36").Value = '38.91'
$ws.Range("D37").Value = '56.27'
$ws.Range("D38").Value = '0.0482'
$ws.Range("D39").Value = '0.0₃0703'
$ws.Range("D40").Value = '0.145'
$ws.Range("D41").Value = '1.00'
$ws.Range("D42").Value = '2.70'
$ws.Range("D43").Value = '2.93'
$ws.Range("D44").Value = '146.08'
$ws.Range("D45").Value = '3.19'
$ws.Range("D46").Value = '4.24'
$ws.Range("D47").Value = '0.302'
$ws.Range("D48").Value = '1.93'
$ws.Range("D49").Value = '2.23'
$ws.Range("D50").Value = '115.23'

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = '  +5.91%  '
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  -1.96%  '
$ws.Range("E7").Value = '  +2.61%  '
$ws.Range("E8").Value = '  +2.36%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  +5.21%  '
$ws.Range("E11").Value = '  +18.63%  '
$ws.Range("E12").Value = '  +39.64%  '
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("E15").Value = '  +2.55%  '
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("E17").Value = '  -2.60%  '
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("E19").Value = '  +1.72%  '
$ws.Range("E20").Value = '  +6.00%  '
$ws.Range("E21").Value = '  -4.97%  '
$ws.Range("E22").Value = '  -5.47%  '
$ws.Range("E23").Value = '  -2.74%  '
$ws.Range("E24").Value = '  -5.22%  '
$ws.Range("E25").Value = '  -3.96%  '
$ws.Range("E26").Value = '  +1.93%  '
$ws.Range("E27").Value = '  -6.84%  '
$ws.Range("E28").Value = '  +1.90%  '
$ws.Range("E29").Value = '  +1.15%  '
$ws.Range("E30").Value = '  +4.00%  '
$ws.Range("E31").Value = '  +1.54%  '
$ws.Range("E33").Value = '  -5.40%  '
$ws.Range("E34").Value = '  -5.10%  '
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E36").Value = '  -4.79%  '
$ws.Range("E37").Value = '  -3.83%  '
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("E39").Value = '  +26.03%  '
$ws.Range("E40").Value = '  +8.63%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("E43").Value = '  -4.61%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("E45").Value = '  -4.49%  '
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("E47").Value = '  -5.80%  '
$ws.Range("E48").Value = '  -6.92%  '
$ws.Range("E49").Value = '  -7.51%  '
$ws.Range("E50").Value = '  +5.58%  '
$ws.Range("E51").Value = '  +8.58%  '

# --- Coin / Link swaps (rows re-ordered by rank in the source feed) ---
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("B45").Value = 'LidoDAOToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
